# Scheduled-runner style market-data refresh: updates the computed
# price/profit columns (H:N) for a handful of leve rows across each
# Disciple-of-the-Hand sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5741.1763
$ws.Range("I51").Value = 1233.6666
$ws.Range("J51").Value = 6707.0713
$ws.Range("K51").Value = 1233.6666
$ws.Range("L51").Value = 6707.0713
$ws.Range("M51").Value = -749.6666
$ws.Range("N51").Value = -7675.0713
$ws.Range("H62").Value = 2530.6956
$ws.Range("I62").Value = 1810.0714
$ws.Range("J62").Value = 3651.6667
$ws.Range("K62").Value = 1810.0714
$ws.Range("L62").Value = 3651.6667
$ws.Range("M62").Value = -1186.0714
$ws.Range("N62").Value = -4899.6667
$ws.Range("H65").Value = 2530.6956
$ws.Range("I65").Value = 1810.0714
$ws.Range("J65").Value = 3651.6667
$ws.Range("K65").Value = 9050.357
$ws.Range("L65").Value = 18258.3335
$ws.Range("M65").Value = -5930.357
$ws.Range("N65").Value = -24498.3335
$ws.Range("H76").Value = 3379.16
$ws.Range("I76").Value = 2998.9268
$ws.Range("K76").Value = 2998.9268
$ws.Range("M76").Value = -2683.9268
$ws.Range("H79").Value = 3379.16
$ws.Range("I79").Value = 2998.9268
$ws.Range("K79").Value = 2998.9268
$ws.Range("M79").Value = -1906.9268
$ws.Range("H132").Value = 38041.035
$ws.Range("I132").Value = 41024.32
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 123072.96
$ws.Range("L132").Value = 2250
$ws.Range("M132").Value = -120542.96
$ws.Range("N132").Value = -7310
$ws.Range("H137").Value = 1268.6511
$ws.Range("I137").Value = 967.3158
$ws.Range("K137").Value = 2901.9474
$ws.Range("M137").Value = -351.9474
$ws.Range("H138").Value = 1594.2174
$ws.Range("I138").Value = 1173.9269
$ws.Range("K138").Value = 3521.7807
$ws.Range("M138").Value = 1618.2193
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3406.4119
$ws.Range("I32").Value = 2499.1392
$ws.Range("J32").Value = 15352.167
$ws.Range("K32").Value = 2499.1392
$ws.Range("L32").Value = 15352.167
$ws.Range("M32").Value = -2212.1392
$ws.Range("N32").Value = -15926.167
$ws.Range("H74").Value = 1209.4667
$ws.Range("I74").Value = 1266
$ws.Range("K74").Value = 1266
$ws.Range("M74").Value = -392
$ws.Range("H77").Value = 1209.4667
$ws.Range("I77").Value = 1266
$ws.Range("K77").Value = 6330
$ws.Range("M77").Value = -1962
$ws.Range("H122").Value = 2793.7693
$ws.Range("I122").Value = 2017
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 6051
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -3601
$ws.Range("N122").Value = -16000
$ws.Range("H132").Value = 1911.3636
$ws.Range("I132").Value = 1089.3871
$ws.Range("J132").Value = 2973.0833
$ws.Range("K132").Value = 3268.1613
$ws.Range("L132").Value = 8919.249899999999
$ws.Range("M132").Value = -738.1612999999998
$ws.Range("N132").Value = -13979.2499
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2441.4614
$ws.Range("I134").Value = 1416.8572
$ws.Range("J134").Value = 3636.8333
$ws.Range("K134").Value = 4250.571599999999
$ws.Range("L134").Value = 10910.4999
$ws.Range("M134").Value = -1715.571599999999
$ws.Range("N134").Value = -15980.4999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 262.8
$ws.Range("I22").Value = 294.66666
$ws.Range("J22").Value = 215
$ws.Range("K22").Value = 294.66666
$ws.Range("L22").Value = 215
$ws.Range("M22").Value = 55.33334000000002
$ws.Range("N22").Value = -915
$ws.Range("H31").Value = 2358.7708
$ws.Range("I31").Value = 2022.8077
$ws.Range("K31").Value = 2022.8077
$ws.Range("M31").Value = -1727.8077
$ws.Range("H34").Value = 2358.7708
$ws.Range("I34").Value = 2022.8077
$ws.Range("K34").Value = 2022.8077
$ws.Range("M34").Value = -1820.8077
$ws.Range("H132").Value = 2011.0714
$ws.Range("I132").Value = 821.2105
$ws.Range("K132").Value = 2463.6315
$ws.Range("M132").Value = 66.36850000000004
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2202
$ws.Range("I5").Value = 902.3333
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 2706.9999
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = -2594.9999
$ws.Range("N5").Value = -30224
$ws.Range("H46").Value = 875.3333
$ws.Range("I46").Value = 125
$ws.Range("J46").Value = 990.7692
$ws.Range("K46").Value = 375
$ws.Range("L46").Value = 2972.3076
$ws.Range("M46").Value = -284
$ws.Range("N46").Value = -3154.3076
$ws.Range("H70").Value = 2260
$ws.Range("I70").Value = 766.6667
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 2300.0001
$ws.Range("L70").Value = 13500
$ws.Range("M70").Value = -1985.0001
$ws.Range("N70").Value = -14130
$ws.Range("H73").Value = 2260
$ws.Range("I73").Value = 766.6667
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 2300.0001
$ws.Range("L73").Value = 13500
$ws.Range("M73").Value = -1208.0001
$ws.Range("N73").Value = -15684
$ws.Range("H135").Value = 2202
$ws.Range("I135").Value = 902.3333
$ws.Range("J135").Value = 10000
$ws.Range("K135").Value = 8120.9997
$ws.Range("L135").Value = 90000
$ws.Range("M135").Value = -5585.9997
$ws.Range("N135").Value = -95070
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.38461
$ws.Range("I2").Value = 58.75
$ws.Range("J2").Value = 146.33333
$ws.Range("K2").Value = 58.75
$ws.Range("L2").Value = 146.33333
$ws.Range("M2").Value = 54.25
$ws.Range("N2").Value = -372.33333
$ws.Range("H132").Value = 1903.9778
$ws.Range("I132").Value = 1658.5555
$ws.Range("J132").Value = 2885.6667
$ws.Range("K132").Value = 4975.666499999999
$ws.Range("L132").Value = 8657.000100000001
$ws.Range("M132").Value = -2445.666499999999
$ws.Range("N132").Value = -13717.0001
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3294.7
$ws.Range("I122").Value = 2234.15
$ws.Range("J122").Value = 4355.25
$ws.Range("K122").Value = 6702.450000000001
$ws.Range("L122").Value = 13065.75
$ws.Range("M122").Value = -4252.450000000001
$ws.Range("N122").Value = -17965.75
$ws.Range("H132").Value = 5467.9395
$ws.Range("I132").Value = 5785.0435
$ws.Range("J132").Value = 4738.6
$ws.Range("K132").Value = 17355.1305
$ws.Range("L132").Value = 14215.8
$ws.Range("M132").Value = -14825.1305
$ws.Range("N132").Value = -19275.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4445.5386
$ws.Range("I107").Value = 3890.182
$ws.Range("J107").Value = 7500
$ws.Range("K107").Value = 11670.546
$ws.Range("L107").Value = 22500
$ws.Range("M107").Value = -9750.545999999998
$ws.Range("N107").Value = -26340
$ws.Range("H122").Value = 29039.945
$ws.Range("I122").Value = 40154.81
$ws.Range("K122").Value = 120464.43
$ws.Range("M122").Value = -118014.43
$ws.Range("H136").Value = 2542.7666
$ws.Range("I136").Value = 963.4
$ws.Range("J136").Value = 5701.5
$ws.Range("K136").Value = 2890.2
$ws.Range("L136").Value = 17104.5
$ws.Range("M136").Value = -340.1999999999998
$ws.Range("N136").Value = -22204.5
$ws.Range("H137").Value = 48083
$ws.Range("J137").Value = 48083
$ws.Range("L137").Value = 48083
$ws.Range("N137").Value = -58283
$ws.Range("H139").Value = 60642.855
$ws.Range("J139").Value = 60642.855
$ws.Range("L139").Value = 60642.855
$ws.Range("N139").Value = -70922.85500000001
